$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This weekly refresh adds two new price observations to the "Perejil"
# series. They land as two separate single-row insertions (not a
# contiguous 2-row block), so everything below each insertion point
# shifts down by one at a time.

# 1) Insert a new row at 17 (old row 17 -> row 18, old row 18 -> row 19, ...).
$ws.Rows.Item(17).Insert()

# Populate the new row 17 by cloning the now-adjacent row 18 (constant/
# static columns A,B,C,E,F,G,H,I,N,O,Q,R) and then set this week's own
# figures (D,J,K,L,M,P).
$ws.Range("A17:R17").Value2 = $ws.Range("A18:R18").Value2
$ws.Range("D17").Value2 = 45117
$ws.Range("J17").Value2 = 56
$ws.Range("K17").Value2 = 3000
$ws.Range("L17").Value2 = 3000
$ws.Range("M17").Value2 = 3000
$ws.Range("P17").Value2 = 1000

# 2) Insert a second new row at 19 (old row 18, now at row 19, -> row 20, ...).
$ws.Rows.Item(19).Insert()

# Populate the new row 19 by cloning the now-adjacent row 20 and then set
# this week's own figures.
$ws.Range("A19:R19").Value2 = $ws.Range("A20:R20").Value2
$ws.Range("D19").Value2 = 45118
$ws.Range("J19").Value2 = 67
$ws.Range("K19").Value2 = 3000
$ws.Range("L19").Value2 = 3000
$ws.Range("M19").Value2 = 3000
$ws.Range("P19").Value2 = 1000
